$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "-"
$ws.Range("D2").Value = "-"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"

# Row 4
$ws.Range("C4").Value = "MEC-3B-Cont. Lóg. Prog. CLP"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "MCT-3A-Eletropneumática"

# Row 6
$ws.Range("C6").Value = "MEC-3B-Cont. Lóg. Prog. CLP"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "MCT-3A-Eletropneumática"

# Row 7
$ws.Range("C7").Value = "MEC-3B-Cont. Lóg. Prog. CLP"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "MCT-3A-Eletropneumática"

# Row 8
$ws.Range("C8").Value = "MEC-3B-Cont. Lóg. Prog. CLP"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "MCT-3A-Eletropneumática"

# Row 11
$ws.Range("B11").Value = "MEC-3A-Comam. Pneumáticos"

# Row 12
$ws.Range("B12").Value = "MEC-3A-Comam. Pneumáticos"
$ws.Range("C12").Value = "-"

# Row 14
$ws.Range("B14").Value = "MEC-3A-Comam. Pneumáticos"
$ws.Range("C14").Value = "-"

# Row 15
$ws.Range("B15").Value = "MEC-3A-Comam. Pneumáticos"
$ws.Range("C15").Value = "-"

# Row 16
$ws.Range("C16").Value = "-"
